$neo4jFile = @'
TC05_CDS_Filter_InstrumentModel-Illumina HiSeq 2000_Neo4jData.xlsx
'@
$webFile = @'
TC05_CDS_Filter_InstrumentModel-Illumina HiSeq 2000_WebData.xlsx
'@
$participantQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2000']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY `Participant ID`LIMIT 100
'@
$sampleQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2000']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@
$fileQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2000']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN 
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER By f.file_name LIMIT 100
'@
$matchFileQuery = @'
MATCH (f:file)
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2000']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,f, s, collect(distinct samp.sample_id) as samp
RETURN
count(distinct s) AS Studies,
count(distinct p) AS Participants,
count(distinct samp) AS Samples,
count(distinct f) AS Files
'@

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The order in which *new* (previously unseen) string values are assigned to
# cells determines the order they are appended to the shared-strings table
# once the now-unreferenced old strings are pruned. To reproduce the target
# table layout (filenames first, then the three short queries, then the
# long "MATCH (f:file)" query last) we assign in that exact order, and for
# each distinct string we touch every cell that references it so the old
# value becomes fully dereferenced.

# 1) Neo4j filename -> D2, D3, D4
$ws.Cells.Item(2,4).Value = $neo4jFile
$ws.Cells.Item(3,4).Value = $neo4jFile
$ws.Cells.Item(4,4).Value = $neo4jFile

# 2) WebData filename -> E2, E3, E4
$ws.Cells.Item(2,5).Value = $webFile
$ws.Cells.Item(3,5).Value = $webFile
$ws.Cells.Item(4,5).Value = $webFile

# 3) Participant query -> B2
$ws.Cells.Item(2,2).Value = $participantQuery

# 4) Sample query -> B3
$ws.Cells.Item(3,2).Value = $sampleQuery

# 5) File query -> B4
$ws.Cells.Item(4,2).Value = $fileQuery

# 6) The big "MATCH (f:file)" stats query -> C2, C3, C4
$ws.Cells.Item(2,3).Value = $matchFileQuery
$ws.Cells.Item(3,3).Value = $matchFileQuery
$ws.Cells.Item(4,3).Value = $matchFileQuery

# Columns D and E were best-fit to their (now longer) contents.
$ws.Columns.Item(4).ColumnWidth = 94
$ws.Columns.Item(5).ColumnWidth = 92.33333333333333

# Selection moved to D3.
$ws.Range("D3").Select()

Write-Host "edit complete"
